$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.542.42'
$ws.Range("E2").Value = '  -0.10%  '

$ws.Range("D3").Value = '1.913.36'
$ws.Range("E3").Value = '  -0.27%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.09'
$ws.Range("E5").Value = '  -0.77%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.11%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4836'
$ws.Range("E7").Value = '  +2.16%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2897'
$ws.Range("E8").Value = '  +0.48%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06871'
$ws.Range("E9").Value = '  +0.73%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '110.78'
$ws.Range("E10").Value = '  +5.68%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.27'
$ws.Range("E11").Value = '  +5.36%  '

$ws.Range("D12").Value = '1.912.15'
$ws.Range("E12").Value = '  -0.38%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07571'
$ws.Range("E13").Value = '  -1.54%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.363'
$ws.Range("E14").Value = '  +2.00%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6692'
$ws.Range("E15").Value = '  +0.17%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '293.08'
$ws.Range("E16").Value = '  +0.64%  '

$ws.Range("D17").Value = '30.548.56'
$ws.Range("E17").Value = '  -0.12%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.01'
$ws.Range("E18").Value = '  +0.79%  '

$ws.Range("E19").Value = '  +0.18%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007589'
$ws.Range("E20").Value = '  +0.12%  '

$ws.Range("D21").Value = '2.170.27'
$ws.Range("E21").Value = '  +0.22%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.499'
$ws.Range("E22").Value = '  -0.16%  '

$ws.Range("E23").Value = '  +0.00%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.395'
$ws.Range("E24").Value = '  +0.80%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.448'
$ws.Range("E25").Value = '  +0.74%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.64'
$ws.Range("E26").Value = '  -1.89%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.21'
$ws.Range("E27").Value = '  -3.75%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.090'
$ws.Range("E28").Value = '  -0.96%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1068'
$ws.Range("E29").Value = '  +0.68%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.437'
$ws.Range("E30").Value = '  +2.90%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.128'
$ws.Range("E31").Value = '  -0.70%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.043'
$ws.Range("E32").Value = '  -0.28%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04974'
$ws.Range("E33").Value = '  -1.04%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7340'
$ws.Range("E34").Value = '  +0.08%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.132'
$ws.Range("E35").Value = '  -0.70%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9999'
$ws.Range("E36").Value = '  +0.06%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.710'
$ws.Range("E37").Value = '  -1.25%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02023'
$ws.Range("E38").Value = '  -1.49%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.685'
$ws.Range("E39").Value = '  -0.11%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.012'
$ws.Range("E40").Value = '  -1.64%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '109.40'
$ws.Range("E41").Value = '  -1.20%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4427'
$ws.Range("E42").Value = '  +1.08%  '

$ws.Range("E43").Value = '  -1.43%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.821'
$ws.Range("E44").Value = '  -0.60%  '

$ws.Range("E45").Value = '  +0.10%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '69.14'
$ws.Range("E46").Value = '  +3.29%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.189'
$ws.Range("E47").Value = '  -0.65%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.238'
$ws.Range("E48").Value = '  -0.73%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '48.00'
$ws.Range("E49").Value = '  +0.40%  '

$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.1224'
$ws.Range("E50").Value = '  -0.12%  '

$ws.Range("B51").Value = 'WOONetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.2529'
$ws.Range("E51").Value = '  +3.25%  '
